# For each numbered group (1-9), the "Competency" and "Discomfort" row
# labels in column A were swapped with each other. The numeric values in
# column B stay in place; only the A-column text labels move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($n = 1; $n -le 9; $n++) {
    $competencyRow = ($n - 1) * 3 + 2
    $discomfortRow = $competencyRow + 1

    $competencyCell = $ws.Cells.Item($competencyRow, 1)
    $discomfortCell = $ws.Cells.Item($discomfortRow, 1)

    $tmp = $competencyCell.Value2
    $competencyCell.Value = $discomfortCell.Value2
    $discomfortCell.Value = $tmp
}
